$d = $word.ActiveDocument
Write-Output $d.Sections.Count
